$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeAddress, $value) {
    $rng = $ws.Range($rangeAddress)
    $rng.NumberFormat = "@"
    $rng.Value = $value
}

# Row 2 - Bitcoin
Set-TextValue "D2" "38.817.34"
Set-TextValue "E2" "  +2.83%  "

# Row 3 - Ethereum
Set-TextValue "D3" "2.105.63"
Set-TextValue "E3" "  +3.27%  "

# Row 4 - TetherUSD
Set-TextValue "E4" "  +0.02%  "

# Row 5 - BNB
Set-TextValue "D5" "228.46"
Set-TextValue "E5" "  +0.60%  "

# Row 6 - XRP
Set-TextValue "E6" "  +2.08%  "

# Row 7 - Solana
Set-TextValue "D7" "60.49"
Set-TextValue "E7" "  +1.55%  "

# Row 9 - Cardano
Set-TextValue "E9" "  +2.18%  "

# Row 10 - Dogecoin
Set-TextValue "D10" "0.0837"
Set-TextValue "E10" "  +0.40%  "

# Row 11 - TRON
Set-TextValue "E11" "  -0.23%  "

# Row 12 - WrappedliquidstakedEther2.0
Set-TextValue "D12" "2.419.29"
Set-TextValue "E12" "  +3.42%  "

# Row 13 - Chainlink
Set-TextValue "D13" "15.00"
Set-TextValue "E13" "  +3.91%  "

# Row 14 - Avalanche
Set-TextValue "D14" "22.28"
Set-TextValue "E14" "  +5.93%  "

# Row 15 - Polygon
Set-TextValue "D15" "0.796"
Set-TextValue "E15" "  +3.06%  "

# Row 16 - Polkadot
Set-TextValue "E16" "  +0.33%  "

# Row 17 - WrappedEther
Set-TextValue "D17" "2.098.85"
Set-TextValue "E17" "  +3.39%  "

# Row 18 - WrappedBTC
Set-TextValue "D18" "38.729.56"
Set-TextValue "E18" "  +2.83%  "

# Row 19 - Litecoin
Set-TextValue "D19" "71.90"
Set-TextValue "E19" "  +3.68%  "

# Row 20 - Uniswap
Set-TextValue "D20" "6.05"
Set-TextValue "E20" "  +1.86%  "

# Row 21 - ShibaInu
Set-TextValue "D21" "0.0₃0836"
Set-TextValue "E21" "  +1.60%  "

# Row 22 - BitcoinCash
Set-TextValue "D22" "226.07"
Set-TextValue "E22" "  +1.01%  "

# Row 23 - Dai
Set-TextValue "E23" "  -0.16%  "

# Row 24 - Toncoin
Set-TextValue "E24" "  -0.25%  "

# Row 25 - PancakeSwap
Set-TextValue "E25" "  +2.75%  "

# Row 26 - Monero
Set-TextValue "D26" "170.63"
Set-TextValue "E26" "  +1.56%  "

# Row 27 - Cosmos
Set-TextValue "D27" "9.50"
Set-TextValue "E27" "  +1.34%  "

# Row 28 - Kaspa
Set-TextValue "E28" "  +6.03%  "

# Row 29 - ImmutableX
Set-TextValue "D29" "1.39"
Set-TextValue "E29" "  +9.37%  "

# Row 30 - EthereumClassic
Set-TextValue "D30" "19.21"
Set-TextValue "E30" "  +2.29%  "

# Row 31 - Stellar
Set-TextValue "E31" "  +0.34%  "

# Row 32 - WEMIXToken
Set-TextValue "E32" "  +4.84%  "

# Row 33 - InternetComputer(DFINITY)
Set-TextValue "E33" "  +6.69%  "

# Row 34 - Filecoin
Set-TextValue "E34" "  +2.95%  "

# Row 35 - Hedera
Set-TextValue "D35" "0.0614"
Set-TextValue "E35" "  +0.93%  "

# Row 36 - LidoDAOToken
Set-TextValue "E36" "  +1.97%  "

# Row 37 - THORChain
Set-TextValue "D37" "6.42"
Set-TextValue "E37" "  -1.92%  "

# Row 38 - RenderToken
Set-TextValue "E38" "  +3.53%  "

# Row 39 - BinanceUSD
Set-TextValue "E39" "  +0.17%  "

# Row 40 - InjectiveProtocol
Set-TextValue "D40" "18.40"
Set-TextValue "E40" "  +1.77%  "

# Row 41 - Maker
Set-TextValue "D41" "1.545.31"
Set-TextValue "E41" "  +0.96%  "

# Row 42 - Aave
Set-TextValue "D42" "101.53"
Set-TextValue "E42" "  +4.63%  "

# Row 43 - VeChain
Set-TextValue "E43" "  +3.41%  "

# Row 44 - was HuobiToken, now Cronos
Set-TextValue "B44" "Cronos"
Set-TextValue "C44" "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue "D44" "0.0927"
Set-TextValue "E44" "  +2.21%  "

# Row 45 - was Cronos, now HuobiToken
Set-TextValue "B45" "HuobiToken"
Set-TextValue "C45" "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue "D45" "2.82"
Set-TextValue "E45" "  -0.77%  "

# Row 46 - FraxShare
Set-TextValue "D46" "7.68"
Set-TextValue "E46" "  +9.01%  "

# Row 47 - FTXToken
Set-TextValue "E47" "  -3.74%  "

# Row 48 - TrustWalletToken
Set-TextValue "D48" "1.12"
Set-TextValue "E48" "  +0.79%  "

# Row 49 - ARBITRUM
Set-TextValue "D49" "1.04"
Set-TextValue "E49" "  +2.77%  "

# Row 50 - MXToken
Set-TextValue "E50" "  +1.93%  "

# Row 51 - RocketPoolETH
Set-TextValue "D51" "2.305.34"
Set-TextValue "E51" "  +3.45%  "
